# Updated cryptos list on Fri Dec  1 13:36:49 UTC 2023 with GitHub Actions
#
# Refresh the Price (column D) and Volume(1h) (column E) columns of the
# cryptos sheet with newly-scraped values. NumberFormat is forced to Text
# on each touched cell immediately before the write so Excel does not
# reinterpret numeric-looking strings (e.g. "22.40") as numbers and drop
# significant trailing digits - the source data must stay literal text,
# matching how it was originally authored.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "38.437.76"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.089.81"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.31"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.77%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.85"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.381"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.35%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.03%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.400.51"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.86"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.40"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +6.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.787"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.16%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.083.45"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "38.334.45"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.37"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.56%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.27%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "225.54"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.43"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.96"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.44"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.136"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.04"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.19%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +8.43%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.82"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +6.81%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0607"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.43"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.82%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.54"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.33%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.45"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.539.15"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.01"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.78%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.35%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.06%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.73"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +8.65%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.27%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.88%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.52%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.287.15"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.19%  "
